$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20:128 down to 21:129
$ws.Rows.Item(20).EntireRow.Insert()

# Populate the newly inserted row 20 with the new record
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 44462
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 100112005
$ws.Range("G20").Value = "Puerro"
$ws.Range("H20").Value = "Azul de Maquehue"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 8000
$ws.Range("N20").Value = "$/docena de paquetes"
$ws.Range("O20").Value = "Provincia de Cautín"
$ws.Range("P20").Value = 667
$ws.Range("Q20").Value = 12
$ws.Range("R20").Value = "Hortaliza"
